# Bump the "Förändrad" (Changed) date in column C by one day for every
# data row (rows 2-496). In the source file every cell in C2:C496 held
# the same serial date value 45181 and the commit changes all of them
# to 45182 (i.e. one day later), leaving everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 496 }

$range = $ws.Range("C2:C" + $lastRow)
$range.Value = 45182
